$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.599.87"
$ws.Range("E2").Value = "  +6.33%  "
$ws.Range("D3").Value = "2.638.18"
$ws.Range("E3").Value = "  +9.96%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'512.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").Value = "'158.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("D7").Value = "'0.993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "2.680.54"
$ws.Range("E9").Value = "  +10.74%  "
$ws.Range("D10").Value = "'6.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").Value = "  +5.99%  "
$ws.Range("E12").Value = "  +4.40%  "
$ws.Range("D14").Value = "3.115.07"
$ws.Range("E14").Value = "  +10.18%  "
$ws.Range("D15").Value = "60.685.78"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").Value = "'21.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.41%  "
$ws.Range("E17").Value = "  +6.38%  "
$ws.Range("D18").Value = "2.678.11"
$ws.Range("E18").Value = "  +10.62%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "'348.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.78%  "
$ws.Range("E21").Value = "  +6.83%  "
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'60.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.78%  "
$ws.Range("E25").Value = "  +5.25%  "
$ws.Range("D26").Value = "2.795.08"
$ws.Range("E26").Value = "  +11.05%  "
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("D28").Value = "'0.992"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "0.0₃0877"
$ws.Range("E29").Value = "  +12.77%  "
$ws.Range("D30").Value = "'7.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +6.02%  "
$ws.Range("D33").Value = "'157.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("D35").Value = "'5.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.77%  "
$ws.Range("E36").Value = "  +10.76%  "
$ws.Range("E37").Value = "  +6.10%  "
$ws.Range("D38").Value = "'313.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.98%  "
$ws.Range("E39").Value = "  +10.50%  "
$ws.Range("D40").Value = "'0.862"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("D41").Value = "'0.848"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +33.59%  "
$ws.Range("E42").Value = "  +7.67%  "
$ws.Range("D43").Value = "'35.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("E44").Value = "  +8.78%  "
$ws.Range("D45").Value = "'0.0584"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.61%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'20.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.42%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'4.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.64%  "
$ws.Range("D50").Value = "2.079.67"
$ws.Range("E50").Value = "  +11.20%  "
$ws.Range("E51").Value = "  +3.97%  "
